$wb = $excel.ActiveWorkbook

# xlPasteFormats: paste only the cell formatting, leaving the destination's
# current value untouched.
$xlPasteFormats = -4122

$passenger = $wb.Worksheets.Item("BCbVT-passenger")
$freight   = $wb.Worksheets.Item("BCbVT-freight")

foreach ($ws in @($passenger, $freight)) {
    # Give the two brand-new header cells (G1's current formatting carries
    # bold + right alignment, same as every other header cell) the same
    # look as the rest of row 1 before touching any cell values.
    $ws.Cells.Item(1, 7).Copy()
    $ws.Cells.Item(1, 8).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item(1, 9).PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false

    # Shift the existing trailing "nonroad vehicle" header out to the new
    # last column, then introduce the two new vehicle types in its place.
    $ws.Cells.Item(1, 9).Value = $ws.Cells.Item(1, 7).Value2
    $ws.Cells.Item(1, 7).Value = "LPG vehicle"
    $ws.Cells.Item(1, 8).Value = "hydrogen vehicle"

    # New data columns default to 0, same as the other vehicle-type columns.
    for ($r = 2; $r -le 7; $r++) {
        $ws.Cells.Item($r, 8).Value = 0
        $ws.Cells.Item($r, 9).Value = 0
    }

    # Label + wrap the corner cell and grow the header row to fit it.
    $ws.Cells.Item(1, 1).Value = "Battery Capacity (MW*hr/vehicle"
    $ws.Cells.Item(1, 1).WrapText = $true
    $ws.Rows.Item(1).RowHeight = 57
}
